$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new "ethanol" unit process rows (55 and 56) ---
# Values are entered in the same order the original author appears to have
# used (column by column across both new rows) so that the shared-string
# table gets populated in the same order as the source edit.

$ws.Range("A55").Value = "eth_box"
$ws.Range("A56").Value = "eth_stoich"

$ws.Range("B55").Value = "chemicals"
$ws.Range("B56").Value = "chemicals"

$ws.Range("D55").Value = "C2H6O"
$ws.Range("D56").Value = "C2H6O"

$ws.Range("C55").Value = "Ethanol (Black Box)"
$ws.Range("C56").Value = "Ethanol (Stoichiometric)"

$ws.Range("F55").Value = "data/chemicals/chemvar.xlsx"
$ws.Range("F56").Value = "data/chemicals/chemvar.xlsx"

$ws.Range("G56").Value = "eth-stoich"
$ws.Range("I56").Value = "eth-stoich"

$ws.Range("G55").Value = "eth-box"
$ws.Range("I55").Value = "eth-box"

$ws.Range("H55").Value = "data/chemicals/chemcals.xlsx"
$ws.Range("H56").Value = "data/chemicals/chemcals.xlsx"

$ws.Range("E55").Value = "outflow"
$ws.Range("E56").Value = "outflow"

# --- Apply the same text number format used by the rest of the table ---
# (Columns B and D are left in the default/general format, matching the
# pattern already used by every other product row in the sheet.)
$ws.Range("A55:A56").NumberFormat = "@"
$ws.Range("C55:C56").NumberFormat = "@"
$ws.Range("E55:E56").NumberFormat = "@"
$ws.Range("F55:F56").NumberFormat = "@"
$ws.Range("G55:G56").NumberFormat = "@"
$ws.Range("H55:H56").NumberFormat = "@"
$ws.Range("I55:I56").NumberFormat = "@"

# --- Update the view state to reflect the author's final selection ---
$ws.Range("H54").Select()
